$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-22 23:48:40'
$ws.Range('O2').Value = '5.8 °C'
$ws.Range('E3').Value = '2026-02-22 23:48:42'
$ws.Range('N3').Value = '0.2 °C 23:11 TU'
$ws.Range('E4').Value = '2026-02-22 23:48:45'
$ws.Range('O4').Value = '11.9 °C'
$ws.Range('E5').Value = '2026-02-22 23:48:47'
$ws.Range('L5').Value = '20.9 km/h - 311º 23:24 TU'
$ws.Range('O5').Value = '5.7 °C'
$ws.Range('E6').Value = '2026-02-22 23:48:50'
$ws.Range('H6').Formula = '="63%"'
$ws.Range('H6').Copy()
$ws.Range('H6').PasteSpecial(-4163)
$ws.Range('E7').Value = '2026-02-22 23:48:53'
$ws.Range('O7').Value = '14.2 °C'
$ws.Range('E8').Value = '2026-02-22 23:48:56'
$ws.Range('J8').Value = '1026.8 hPa'
$ws.Range('E9').Value = '2026-02-22 23:48:58'
$ws.Range('O9').Value = '10.7 °C'
$ws.Range('E10').Value = '2026-02-22 23:49:01'
$ws.Range('O10').Value = '9.6 °C'
$ws.Range('E11').Value = '2026-02-22 23:49:03'
$ws.Range('H11').Formula = '="67%"'
$ws.Range('H11').Copy()
$ws.Range('H11').PasteSpecial(-4163)
$ws.Range('O11').Value = '8.4 °C'
$ws.Range('E12').Value = '2026-02-22 23:49:06'
$ws.Range('E13').Value = '2026-02-22 23:49:08'
$ws.Range('H13').Formula = '="63%"'
$ws.Range('H13').Copy()
$ws.Range('H13').PasteSpecial(-4163)
$ws.Range('O13').Value = '6.3 °C'
$ws.Range('E14').Value = '2026-02-22 23:49:11'
$ws.Range('H14').Formula = '="76%"'
$ws.Range('H14').Copy()
$ws.Range('H14').PasteSpecial(-4163)
$ws.Range('E15').Value = '2026-02-22 23:49:14'
$ws.Range('E16').Value = '2026-02-22 23:49:16'
$ws.Range('N16').Value = '2.3 °C 23:10 TU'
$ws.Range('O16').Value = '5.2 °C'
$ws.Range('E17').Value = '2026-02-22 23:49:19'
$ws.Range('N17').Value = '7.1 °C 23:14 TU'
$ws.Range('E18').Value = '2026-02-22 23:49:21'
$ws.Range('H18').Formula = '="79%"'
$ws.Range('H18').Copy()
$ws.Range('H18').PasteSpecial(-4163)
$ws.Range('O18').Value = '9.7 °C'
$ws.Range('E19').Value = '2026-02-22 23:49:24'
$ws.Range('H19').Formula = '="49%"'
$ws.Range('H19').Copy()
$ws.Range('H19').PasteSpecial(-4163)
$ws.Range('O19').Value = '11.8 °C'
$ws.Range('E20').Value = '2026-02-22 23:49:27'
$ws.Range('E21').Value = '2026-02-22 23:49:29'
$ws.Range('O21').Value = '9.0 °C'
$ws.Range('E22').Value = '2026-02-22 23:49:32'
$ws.Range('N22').Value = '1.9 °C 23:28 TU'
$ws.Range('O22').Value = '4.4 °C'
$ws.Range('E23').Value = '2026-02-22 23:49:35'
$ws.Range('E24').Value = '2026-02-22 23:49:37'
$ws.Range('O24').Value = '7.5 °C'
$ws.Range('E25').Value = '2026-02-22 23:49:40'
$ws.Range('O25').Value = '6.8 °C'
$ws.Range('E26').Value = '2026-02-22 23:49:43'
$ws.Range('J26').Value = '1026.3 hPa'
$ws.Range('N26').Value = '6.8 °C 23:18 TU'
$ws.Range('O26').Value = '10.8 °C'
$ws.Range('E27').Value = '2026-02-22 23:49:45'
$ws.Range('O27').Value = '6.5 °C'
$ws.Range('E28').Value = '2026-02-22 23:49:48'
$ws.Range('O28').Value = '10.1 °C'
$ws.Range('E29').Value = '2026-02-22 23:49:51'
$ws.Range('O29').Value = '9.4 °C'
$ws.Range('E30').Value = '2026-02-22 23:49:53'
$ws.Range('H30').Formula = '="75%"'
$ws.Range('H30').Copy()
$ws.Range('H30').PasteSpecial(-4163)
$ws.Range('O30').Value = '11.9 °C'
$ws.Range('E31').Value = '2026-02-22 23:49:56'
$ws.Range('L31').Value = '47.5 km/h - 335º 23:09 TU'
$ws.Range('E32').Value = '2026-02-22 23:49:59'
$ws.Range('H32').Formula = '="74%"'
$ws.Range('H32').Copy()
$ws.Range('H32').PasteSpecial(-4163)
$ws.Range('O32').Value = '5.4 °C'
$ws.Range('E33').Value = '2026-02-22 23:50:01'
$ws.Range('E34').Value = '2026-02-22 23:50:04'
$ws.Range('E35').Value = '2026-02-22 23:50:07'
$ws.Range('J35').Value = '1028.4 hPa'
$ws.Range('E36').Value = '2026-02-22 23:50:10'
$ws.Range('J36').Value = '1027.2 hPa'
$ws.Range('O36').Value = '11.4 °C'
$ws.Range('E37').Value = '2026-02-22 23:50:12'
$ws.Range('O37').Value = '7.6 °C'
$ws.Range('E38').Value = '2026-02-22 23:50:15'
$ws.Range('H38').Formula = '="68%"'
$ws.Range('H38').Copy()
$ws.Range('H38').PasteSpecial(-4163)
$ws.Range('O38').Value = '11.3 °C'
$ws.Range('E39').Value = '2026-02-22 23:50:18'
$ws.Range('K39').Value = '16.6 MJ/m2'
$ws.Range('E40').Value = '2026-02-22 23:50:20'
$ws.Range('O40').Value = '9.4 °C'
$ws.Range('E41').Value = '2026-02-22 23:50:22'
$ws.Range('H41').Formula = '="80%"'
$ws.Range('H41').Copy()
$ws.Range('H41').PasteSpecial(-4163)
$ws.Range('J41').Value = '1027.5 hPa'
$ws.Range('O41').Value = '11.0 °C'
$ws.Range('E42').Value = '2026-02-22 23:50:25'
$ws.Range('H42').Formula = '="83%"'
$ws.Range('H42').Copy()
$ws.Range('H42').PasteSpecial(-4163)
$ws.Range('O42').Value = '10.3 °C'
$ws.Range('E43').Value = '2026-02-22 23:50:28'
$ws.Range('O43').Value = '9.2 °C'
$ws.Range('E44').Value = '2026-02-22 23:50:30'
$ws.Range('L44').Value = '37.4 km/h - 75º 23:20 TU'
$ws.Range('E45').Value = '2026-02-22 23:50:33'
$ws.Range('J45').Value = '1029.3 hPa'
$ws.Range('O45').Value = '8.3 °C'
$ws.Range('E46').Value = '2026-02-22 23:50:36'
$ws.Range('H46').Formula = '="77%"'
$ws.Range('H46').Copy()
$ws.Range('H46').PasteSpecial(-4163)
$ws.Range('O46').Value = '8.7 °C'

$excel.CutCopyMode = $false
